$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("handhelds_mtbenson")

# Replace every "needed" placeholder value with "placeholder"
$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    if ($cell.Value() -eq "needed") {
        $cell.Value = "placeholder"
    }
}

# Update the active selection to F9, matching the saved view state
$ws.Range("F9").Select()
